{"js": "// Word JS API (Office.js) edit script.\n// Locate the FIRST occurrence of the \"Sometimes a monster...\" clue paragraph\n// (the document repeats this clue block multiple times; only the first\n// occurrence is modified, per the target diff) and rewrite it plus the\n// nearby \"Use this spell...\" paragraph to match the target content.\n\nconst body = context.document.body;\n\n// --- Find the first occurrence of the target clue block -------------------\nconst searchResults = body.search(\"Sometimes a monster just needs a friend\", { matchCase: false });\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not locate target paragraph 'Sometimes a monster...'\");\n}\n\nconst firstHit = searchResults.items[0];\nconst targetPara = firstHit.paragraphs.getFirst();\n\n// --- Rewrite the \"Sometimes a monster...\" paragraph -----------------------\n// It now ends after \"spell.\" and is followed by a new empty paragraph and a\n// new paragraph with a concrete cp example.\nconst para1Ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>Sometimes monster</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>s</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> just need a friend. Create a copy of the original monster </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>with</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> the </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>cp</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>spell.</w:t></w:r>\n</w:p>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n</w:p>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>An example of this</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>spell</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> is </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>c</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>p file1.txt new_file.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>.</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> This would create a new file named </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>new_file.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> that is an exact copy of </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>file</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>1.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>.</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntargetPara.insertOoxml(para1Ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-locate the \"Use this spell to create a duplicate monster called ...\"\n// paragraph with a fresh search: the earlier `targetPara` reference's\n// siblings are no longer reliable once the surrounding paragraph structure\n// has been rewritten by the insertOoxml(\"Replace\") call above.\nconst useThisSpellResults = body.search(\"Use this spell to create a duplicate monster called\", { matchCase: false });\nawait context.sync();\n\nif (useThisSpellResults.items.length === 0) {\n  throw new Error(\"Could not locate 'Use this spell...' paragraph after first edit.\");\n}\n\nconst useThisSpellPara = useThisSpellResults.items[0].paragraphs.getFirst();\n\n// --- Rewrite the \"Use this spell to create a duplicate ...\" paragraph -----\nconst para2Ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">Use this spell to create a duplicate </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">of </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>monster.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> called </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>monster_friend.txt</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nuseThisSpellPara.insertOoxml(para2Ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Locate the FIRST occurrence of the \"Sometimes a monster...\" clue paragraph\n# (the document repeats this clue block multiple times; only the first\n# occurrence is modified, per the target diff) and rewrite it plus the\n# nearby \"Use this spell...\" paragraph to match the target content.\n\n$d = $word.ActiveDocument\n\n# --- Find the first occurrence of the target clue block --------------------\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"Sometimes a monster just needs a friend\")\nif (-not $found) {\n    throw \"Could not locate target paragraph 'Sometimes a monster...'\"\n}\n$targetParaRange = $findRange.Paragraphs(1).Range\n\n# --- Rewrite the \"Sometimes a monster...\" paragraph ------------------------\n# It now ends after \"spell.\" and is followed by a new empty paragraph and a\n# new paragraph with a concrete cp example.\n$para1Xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>Sometimes monster</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>s</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> just need a friend. Create a copy of the original monster </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>with</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> the </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>cp</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>spell.</w:t></w:r>\n</w:p>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n</w:p>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>An example of this</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>spell</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> is </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>c</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>p file1.txt new_file.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>.</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> This would create a new file named </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>new_file.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> that is an exact copy of </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>file</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>1.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t>.</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$targetParaRange.InsertXML($para1Xml)\n\n# --- Rewrite the \"Use this spell to create a duplicate ...\" paragraph ------\n# Re-locate it with a fresh Find: the earlier range/paragraph handles are no\n# longer reliable once the surrounding paragraph structure has been\n# rewritten by the InsertXML call above.\n$findRange2 = $d.Content\n$findRange2.Find.ClearFormatting()\n$found2 = $findRange2.Find.Execute(\"Use this spell to create a duplicate monster called\")\nif (-not $found2) {\n    throw \"Could not locate 'Use this spell...' paragraph after first edit.\"\n}\n$useThisSpellParaRange = $findRange2.Paragraphs(1).Range\n\n$para2Xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr><w:ind w:left=\"360\"/><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr></w:pPr>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">Use this spell to create a duplicate </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\">of </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>monster.txt</w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Garamond\" w:hAnsi=\"Garamond\"/></w:rPr><w:t xml:space=\"preserve\"> called </w:t></w:r>\n  <w:r><w:rPr><w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>monster_friend.txt</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$useThisSpellParaRange.InsertXML($para2Xml)\n"}
